$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 25642646
$ws.Range("I137").Value = 1073
$ws.Range("J137").Value = 125003740
$ws.Range("K137").Value = 3219
$ws.Range("L137").Value = 375011220
$ws.Range("M137").Value = -669
$ws.Range("N137").Value = -375016320

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6930.3105
$ws.Range("I32").Value = 4019.6365
$ws.Range("J32").Value = 16078.143
$ws.Range("K32").Value = 4019.6365
$ws.Range("L32").Value = 16078.143
$ws.Range("M32").Value = -3732.6365
$ws.Range("N32").Value = -16652.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 969
$ws.Range("I45").Value = 958.8
$ws.Range("J45").Value = 994.5
$ws.Range("K45").Value = 958.8
$ws.Range("L45").Value = 994.5
$ws.Range("M45").Value = -581.8
$ws.Range("N45").Value = -1748.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1638.1613
$ws.Range("I61").Value = 1022.6539
$ws.Range("J61").Value = 4838.8
$ws.Range("K61").Value = 1022.6539
$ws.Range("L61").Value = 4838.8
$ws.Range("M61").Value = -810.6539
$ws.Range("N61").Value = -5262.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2039.2333
$ws.Range("I132").Value = 1864.3043
$ws.Range("J132").Value = 2614
$ws.Range("K132").Value = 5592.9129
$ws.Range("L132").Value = 7842
$ws.Range("M132").Value = -3062.9129
$ws.Range("N132").Value = -12902

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1638.1613
$ws.Range("I136").Value = 1022.6539
$ws.Range("J136").Value = 4838.8
$ws.Range("K136").Value = 3067.9617
$ws.Range("L136").Value = 14516.4
$ws.Range("M136").Value = -517.9616999999998
$ws.Range("N136").Value = -19616.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 17393.334
$ws.Range("I35").Value = 9000
$ws.Range("J35").Value = 21590
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 21590
$ws.Range("M35").Value = -8690
$ws.Range("N35").Value = -22210

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1779
$ws.Range("I107").Value = 1862.8235
$ws.Range("J107").Value = 1304
$ws.Range("K107").Value = 1862.8235
$ws.Range("L107").Value = 1304
$ws.Range("M107").Value = 57.17650000000003
$ws.Range("N107").Value = -5144

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 45772.81
$ws.Range("I134").Value = 55241.855
$ws.Range("J134").Value = 6002.8
$ws.Range("K134").Value = 165725.565
$ws.Range("L134").Value = 18008.4
$ws.Range("M134").Value = -163190.565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 3000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 3000
$ws.Range("N23").Value = -3480
$ws.Range("M23").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 3000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3384
$ws.Range("M27").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 12500
$ws.Range("I75").Value = 5000
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 5000
$ws.Range("L75").Value = 20000
$ws.Range("M75").Value = -4002
$ws.Range("N75").Value = -21996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H78").Value = 12500
$ws.Range("I78").Value = 5000
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 15000
$ws.Range("L78").Value = 60000
$ws.Range("M78").Value = -10008
$ws.Range("N78").Value = -69984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2602.5
$ws.Range("I122").Value = 2421.2
$ws.Range("J122").Value = 2829.125
$ws.Range("K122").Value = 7263.599999999999
$ws.Range("L122").Value = 8487.375
$ws.Range("M122").Value = -4813.599999999999
$ws.Range("N122").Value = -13387.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2571.1428
$ws.Range("I132").Value = 1750.1428
$ws.Range("J132").Value = 4213.143
$ws.Range("K132").Value = 5250.428400000001
$ws.Range("L132").Value = 12639.429
$ws.Range("M132").Value = -2720.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 275.42856
$ws.Range("I61").Value = 242.66667
$ws.Range("J61").Value = 300
$ws.Range("K61").Value = 728.00001
$ws.Range("L61").Value = 900
$ws.Range("M61").Value = -513.00001
$ws.Range("N61").Value = -1330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 27214.223
$ws.Range("I137").Value = 2142.3333
$ws.Range("J137").Value = 36331.273
$ws.Range("K137").Value = 6426.999899999999
$ws.Range("L137").Value = 108993.819
$ws.Range("M137").Value = -1326.999899999999
$ws.Range("N137").Value = -119193.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I21").Value = 1952.381
$ws.Range("J21").Value = 2000333.4
$ws.Range("K21").Value = 1952.381
$ws.Range("L21").Value = 2000333.4
$ws.Range("M21").Value = -1779.381
$ws.Range("N21").Value = -2000679.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I30").Value = 1952.381
$ws.Range("J30").Value = 2000333.4
$ws.Range("K30").Value = 1952.381
$ws.Range("L30").Value = 2000333.4
$ws.Range("M30").Value = -1847.381
$ws.Range("N30").Value = -2000543.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 625.5714
$ws.Range("I107").Value = 394
$ws.Range("J107").Value = 718.2
$ws.Range("K107").Value = 394
$ws.Range("L107").Value = 718.2
$ws.Range("M107").Value = 1526
$ws.Range("N107").Value = -4558.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3255.68
$ws.Range("I132").Value = 3139.647
$ws.Range("J132").Value = 3502.25
$ws.Range("K132").Value = 9418.940999999999
$ws.Range("L132").Value = 10506.75
$ws.Range("M132").Value = -6888.940999999999
$ws.Range("N132").Value = -15566.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 27400
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 27400
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 27400
$ws.Range("N121").Value = -30894

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2375.9285
$ws.Range("I132").Value = 1220.625
$ws.Range("J132").Value = 3916.3333
$ws.Range("K132").Value = 3661.875
$ws.Range("L132").Value = 11748.9999
$ws.Range("M132").Value = -1131.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16665.691
$ws.Range("I136").Value = 17776.75
$ws.Range("J136").Value = 3333
$ws.Range("K136").Value = 53330.25
$ws.Range("L136").Value = 9999
$ws.Range("M136").Value = -50780.25
